$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-27 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-28 Sunday", 2) | Out-Null
$d.Content.Find.Execute("36-7=29", $true, $false, $false, $false, $false, $true, 1, $false, "92-5=87", 2) | Out-Null
$d.Content.Find.Execute("39-23=16", $true, $false, $false, $false, $false, $true, 1, $false, "63-61=2", 2) | Out-Null
$d.Content.Find.Execute("58-50=8", $true, $false, $false, $false, $false, $true, 1, $false, "36+36=72", 2) | Out-Null
$d.Content.Find.Execute("45-44=1", $true, $false, $false, $false, $false, $true, 1, $false, "83-44=39", 2) | Out-Null
$d.Content.Find.Execute("40-12=28", $true, $false, $false, $false, $false, $true, 1, $false, "69-28=41", 2) | Out-Null
$d.Content.Find.Execute("68+28=96", $true, $false, $false, $false, $false, $true, 1, $false, "0-0=0", 2) | Out-Null
$d.Content.Find.Execute("36+57=93", $true, $false, $false, $false, $false, $true, 1, $false, "2+17=19", 2) | Out-Null
$d.Content.Find.Execute("60-8=52", $true, $false, $false, $false, $false, $true, 1, $false, "42+1=43", 2) | Out-Null
$d.Content.Find.Execute("88-65=23", $true, $false, $false, $false, $false, $true, 1, $false, "8+63=71", 2) | Out-Null
$d.Content.Find.Execute("81-79=2", $true, $false, $false, $false, $false, $true, 1, $false, "81-20=61", 2) | Out-Null
$d.Content.Find.Execute("96-94=2", $true, $false, $false, $false, $false, $true, 1, $false, "11+61=72", 2) | Out-Null
$d.Content.Find.Execute("14-8=6", $true, $false, $false, $false, $false, $true, 1, $false, "33+24=57", 2) | Out-Null
$d.Content.Find.Execute("18+2=20", $true, $false, $false, $false, $false, $true, 1, $false, "54+16=70", 2) | Out-Null
$d.Content.Find.Execute("10+81=91", $true, $false, $false, $false, $false, $true, 1, $false, "48-35=13", 2) | Out-Null
$d.Content.Find.Execute("10+87=97", $true, $false, $false, $false, $false, $true, 1, $false, "97-27=70", 2) | Out-Null
$d.Content.Find.Execute("64-5=59", $true, $false, $false, $false, $false, $true, 1, $false, "18+22=40", 2) | Out-Null
$d.Content.Find.Execute("22+26=48", $true, $false, $false, $false, $false, $true, 1, $false, "86-36=50", 2) | Out-Null
$d.Content.Find.Execute("20+45=65", $true, $false, $false, $false, $false, $true, 1, $false, "49-25=24", 2) | Out-Null
$d.Content.Find.Execute("30-2=28", $true, $false, $false, $false, $false, $true, 1, $false, "22+60=82", 2) | Out-Null
$d.Content.Find.Execute("4+73=77", $true, $false, $false, $false, $false, $true, 1, $false, "30+65=95", 2) | Out-Null
$d.Content.Find.Execute("40-20=20", $true, $false, $false, $false, $false, $true, 1, $false, "17+22=39", 2) | Out-Null
$d.Content.Find.Execute("27-27=0", $true, $false, $false, $false, $false, $true, 1, $false, "14+52=66", 2) | Out-Null
$d.Content.Find.Execute("63-52=11", $true, $false, $false, $false, $false, $true, 1, $false, "43+4=47", 2) | Out-Null
$d.Content.Find.Execute("75+17=92", $true, $false, $false, $false, $false, $true, 1, $false, "57-56=1", 2) | Out-Null
$d.Content.Find.Execute("42+8=50", $true, $false, $false, $false, $false, $true, 1, $false, "23+54=77", 2) | Out-Null
$d.Content.Find.Execute("14+9=23", $true, $false, $false, $false, $false, $true, 1, $false, "75-33=42", 2) | Out-Null
$d.Content.Find.Execute("49-18=31", $true, $false, $false, $false, $false, $true, 1, $false, "51-40=11", 2) | Out-Null
$d.Content.Find.Execute("0+85=85", $true, $false, $false, $false, $false, $true, 1, $false, "90+3=93", 2) | Out-Null
$d.Content.Find.Execute("18+6=24", $true, $false, $false, $false, $false, $true, 1, $false, "42+22=64", 2) | Out-Null
$d.Content.Find.Execute("84-35=49", $true, $false, $false, $false, $false, $true, 1, $false, "97-19=78", 2) | Out-Null
$d.Content.Find.Execute("37+13=50", $true, $false, $false, $false, $false, $true, 1, $false, "48-30=18", 2) | Out-Null
$d.Content.Find.Execute("33+33=66", $true, $false, $false, $false, $false, $true, 1, $false, "48-5=43", 2) | Out-Null
$d.Content.Find.Execute("94-29=65", $true, $false, $false, $false, $false, $true, 1, $false, "13+70=83", 2) | Out-Null
$d.Content.Find.Execute("19-9=10", $true, $false, $false, $false, $false, $true, 1, $false, "85-40=45", 2) | Out-Null
$d.Content.Find.Execute("52+23=75", $true, $false, $false, $false, $false, $true, 1, $false, "25-9=16", 2) | Out-Null
$d.Content.Find.Execute("35+0=35", $true, $false, $false, $false, $false, $true, 1, $false, "1+83=84", 2) | Out-Null
$d.Content.Find.Execute("78-49=29", $true, $false, $false, $false, $false, $true, 1, $false, "98-79=19", 2) | Out-Null
$d.Content.Find.Execute("59-6=53", $true, $false, $false, $false, $false, $true, 1, $false, "30+14=44", 2) | Out-Null
$d.Content.Find.Execute("96-34=62", $true, $false, $false, $false, $false, $true, 1, $false, "95-41=54", 2) | Out-Null
$d.Content.Find.Execute("48+50=98", $true, $false, $false, $false, $false, $true, 1, $false, "34+52=86", 2) | Out-Null
$d.Content.Find.Execute("89-74=15", $true, $false, $false, $false, $false, $true, 1, $false, "23+44=67", 2) | Out-Null
$d.Content.Find.Execute("18+49=67", $true, $false, $false, $false, $false, $true, 1, $false, "63+3=66", 2) | Out-Null
$d.Content.Find.Execute("38-36=2", $true, $false, $false, $false, $false, $true, 1, $false, "88-20=68", 2) | Out-Null
$d.Content.Find.Execute("78-44=34", $true, $false, $false, $false, $false, $true, 1, $false, "93-60=33", 2) | Out-Null
$d.Content.Find.Execute("46-32=14", $true, $false, $false, $false, $false, $true, 1, $false, "69+1=70", 2) | Out-Null
$d.Content.Find.Execute("61-35=26", $true, $false, $false, $false, $false, $true, 1, $false, "95-50=45", 2) | Out-Null
$d.Content.Find.Execute("38+60=98", $true, $false, $false, $false, $false, $true, 1, $false, "74-27=47", 2) | Out-Null
$d.Content.Find.Execute("30+20=50", $true, $false, $false, $false, $false, $true, 1, $false, "38-26=12", 2) | Out-Null
$d.Content.Find.Execute("37+28=65", $true, $false, $false, $false, $false, $true, 1, $false, "73-50=23", 2) | Out-Null
$d.Content.Find.Execute("36+8=44", $true, $false, $false, $false, $false, $true, 1, $false, "57-56=1", 2) | Out-Null
$d.Content.Find.Execute("95-27=68", $true, $false, $false, $false, $false, $true, 1, $false, "46-18=28", 2) | Out-Null
$d.Content.Find.Execute("7+24=31", $true, $false, $false, $false, $false, $true, 1, $false, "25+31=56", 2) | Out-Null
$d.Content.Find.Execute("16+72=88", $true, $false, $false, $false, $false, $true, 1, $false, "53-51=2", 2) | Out-Null
$d.Content.Find.Execute("97-45=52", $true, $false, $false, $false, $false, $true, 1, $false, "74+22=96", 2) | Out-Null
$d.Content.Find.Execute("27+37=64", $true, $false, $false, $false, $false, $true, 1, $false, "78-29=49", 2) | Out-Null
$d.Content.Find.Execute("33-3=30", $true, $false, $false, $false, $false, $true, 1, $false, "17+36=53", 2) | Out-Null
$d.Content.Find.Execute("89-84=5", $true, $false, $false, $false, $false, $true, 1, $false, "24+28=52", 2) | Out-Null
$d.Content.Find.Execute("82-0=82", $true, $false, $false, $false, $false, $true, 1, $false, "74-20=54", 2) | Out-Null
$d.Content.Find.Execute("52+8=60", $true, $false, $false, $false, $false, $true, 1, $false, "51-21=30", 2) | Out-Null
$d.Content.Find.Execute("66-61=5", $true, $false, $false, $false, $false, $true, 1, $false, "20+61=81", 2) | Out-Null
$d.Content.Find.Execute("83-60=23", $true, $false, $false, $false, $false, $true, 1, $false, "28-13=15", 2) | Out-Null
$d.Content.Find.Execute("30+60=90", $true, $false, $false, $false, $false, $true, 1, $false, "9+52=61", 2) | Out-Null
$d.Content.Find.Execute("0+1=1", $true, $false, $false, $false, $false, $true, 1, $false, "56-47=9", 2) | Out-Null
$d.Content.Find.Execute("83-16=67", $true, $false, $false, $false, $false, $true, 1, $false, "23+66=89", 2) | Out-Null
$d.Content.Find.Execute("6+12=18", $true, $false, $false, $false, $false, $true, 1, $false, "58-6=52", 2) | Out-Null
$d.Content.Find.Execute("40+48=88", $true, $false, $false, $false, $false, $true, 1, $false, "60-29=31", 2) | Out-Null
$d.Content.Find.Execute("60-33=27", $true, $false, $false, $false, $false, $true, 1, $false, "59+11=70", 2) | Out-Null
$d.Content.Find.Execute("51-34=17", $true, $false, $false, $false, $false, $true, 1, $false, "96-48=48", 2) | Out-Null
$d.Content.Find.Execute("68-32=36", $true, $false, $false, $false, $false, $true, 1, $false, "97-56=41", 2) | Out-Null
$d.Content.Find.Execute("55-8=47", $true, $false, $false, $false, $false, $true, 1, $false, "63-32=31", 2) | Out-Null
$d.Content.Find.Execute("54-40=14", $true, $false, $false, $false, $false, $true, 1, $false, "80+9=89", 2) | Out-Null
$d.Content.Find.Execute("20-14=6", $true, $false, $false, $false, $false, $true, 1, $false, "34+24=58", 2) | Out-Null
$d.Content.Find.Execute("35+62=97", $true, $false, $false, $false, $false, $true, 1, $false, "69-47=22", 2) | Out-Null
$d.Content.Find.Execute("7+75=82", $true, $false, $false, $false, $false, $true, 1, $false, "55+38=93", 2) | Out-Null
$d.Content.Find.Execute("7+41=48", $true, $false, $false, $false, $false, $true, 1, $false, "42+40=82", 2) | Out-Null
$d.Content.Find.Execute("95+4=99", $true, $false, $false, $false, $false, $true, 1, $false, "56-5=51", 2) | Out-Null
$d.Content.Find.Execute("21+17=38", $true, $false, $false, $false, $false, $true, 1, $false, "94-60=34", 2) | Out-Null
$d.Content.Find.Execute("54-24=30", $true, $false, $false, $false, $false, $true, 1, $false, "4+19=23", 2) | Out-Null
$d.Content.Find.Execute("45+21=66", $true, $false, $false, $false, $false, $true, 1, $false, "77-45=32", 2) | Out-Null
$d.Content.Find.Execute("51+16=67", $true, $false, $false, $false, $false, $true, 1, $false, "66+8=74", 2) | Out-Null
$d.Content.Find.Execute("7+35=42", $true, $false, $false, $false, $false, $true, 1, $false, "34+1=35", 2) | Out-Null
$d.Content.Find.Execute("43-15=28", $true, $false, $false, $false, $false, $true, 1, $false, "79+6=85", 2) | Out-Null
$d.Content.Find.Execute("91-58=33", $true, $false, $false, $false, $false, $true, 1, $false, "38-27=11", 2) | Out-Null
$d.Content.Find.Execute("69-55=14", $true, $false, $false, $false, $false, $true, 1, $false, "44-17=27", 2) | Out-Null
$d.Content.Find.Execute("87-71=16", $true, $false, $false, $false, $false, $true, 1, $false, "59+6=65", 2) | Out-Null
$d.Content.Find.Execute("23-20=3", $true, $false, $false, $false, $false, $true, 1, $false, "96+0=96", 2) | Out-Null
$d.Content.Find.Execute("78-66=12", $true, $false, $false, $false, $false, $true, 1, $false, "87-54=33", 2) | Out-Null
$d.Content.Find.Execute("9+37=46", $true, $false, $false, $false, $false, $true, 1, $false, "29+5=34", 2) | Out-Null
$d.Content.Find.Execute("56+25=81", $true, $false, $false, $false, $false, $true, 1, $false, "80-45=35", 2) | Out-Null
$d.Content.Find.Execute("71-41=30", $true, $false, $false, $false, $false, $true, 1, $false, "3+11=14", 2) | Out-Null
$d.Content.Find.Execute("0+74=74", $true, $false, $false, $false, $false, $true, 1, $false, "5+44=49", 2) | Out-Null
$d.Content.Find.Execute("63-47=16", $true, $false, $false, $false, $false, $true, 1, $false, "93-92=1", 2) | Out-Null
$d.Content.Find.Execute("93-11=82", $true, $false, $false, $false, $false, $true, 1, $false, "6+22=28", 2) | Out-Null
$d.Content.Find.Execute("9+61=70", $true, $false, $false, $false, $false, $true, 1, $false, "65+17=82", 2) | Out-Null
$d.Content.Find.Execute("44-18=26", $true, $false, $false, $false, $false, $true, 1, $false, "31+43=74", 2) | Out-Null
$d.Content.Find.Execute("14+8=22", $true, $false, $false, $false, $false, $true, 1, $false, "71-69=2", 2) | Out-Null
$d.Content.Find.Execute("32+60=92", $true, $false, $false, $false, $false, $true, 1, $false, "20+40=60", 2) | Out-Null
$d.Content.Find.Execute("60-2=58", $true, $false, $false, $false, $false, $true, 1, $false, "57+1=58", 2) | Out-Null
$d.Content.Find.Execute("3+58=61", $true, $false, $false, $false, $false, $true, 1, $false, "38+3=41", 2) | Out-Null
$d.Content.Find.Execute("28+46=74", $true, $false, $false, $false, $false, $true, 1, $false, "82-24=58", 2) | Out-Null
